# This workbook holds a weekly-updated price history for "Choclo" at the
# "Terminal Hortofrutícola Agro Chillán" market. A new weekly price record
# (the most recent one, dated 2023-08-03) is inserted above the existing
# top record (dated 2022-09-28) for the same Variedad/Calidad/Origen group,
# pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 236; rows 236-346 shift down to 237-347.
$ws.Rows.Item(236).Insert()

# The row that used to be row 235 (Dulce o Americano / Primera / Arica y
# Parinacota, dated 2022-09-28) now belongs at row 236 - recreate it there.
$ws.Range("A236").Value2 = 7
$ws.Range("B236").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C236").Value2 = "Ñuble"
$ws.Range("D236").Value2 = 44832
$ws.Range("E236").Value2 = 16
$ws.Range("F236").Value2 = 100112024
$ws.Range("G236").Value2 = "Choclo"
$ws.Range("H236").Value2 = "Dulce o Americano"
$ws.Range("I236").Value2 = "Primera"
$ws.Range("J236").Value2 = 50
$ws.Range("K236").Value2 = 30000
$ws.Range("L236").Value2 = 30000
$ws.Range("M236").Value2 = 30000
$ws.Range("N236").Value2 = "`$/malla 70 unidades"
$ws.Range("O236").Value2 = "Región de Arica y Parinacota"
$ws.Range("P236").Value2 = 429
$ws.Range("Q236").Value2 = 70
$ws.Range("R236").Value2 = "Hortaliza"

# Row 235 becomes the brand-new, most recent price record (same
# Variedad/Calidad/Unidad/Origen, newer date and updated prices).
$ws.Range("D235").Value2 = 45141
$ws.Range("J235").Value2 = 60
$ws.Range("K235").Value2 = 44000
$ws.Range("L235").Value2 = 44000
$ws.Range("M235").Value2 = 44000
$ws.Range("P235").Value2 = 629
